$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E/F/H cells for the specified rows, forcing text storage
# (matches original workbook convention where all cells are stored as text,
# even for numeric-looking values like "0.55" or date-like strings)

$ws.Range("E95:H95").NumberFormat = "@"
$ws.Range("E95").Value = "2025-12-31T02:54:20.008-0500"
$ws.Range("F95").Value = "0.55"
$ws.Range("H95").Value = "203.92"
$ws.Range("E95:H95").ClearFormats()

$ws.Range("E97:H97").NumberFormat = "@"
$ws.Range("E97").Value = "2025-12-31T00:06:36.694-0500"
$ws.Range("F97").Value = "6.38"
$ws.Range("H97").Value = "206.46"
$ws.Range("E97:H97").ClearFormats()

$ws.Range("E98:H98").NumberFormat = "@"
$ws.Range("E98").Value = "2025-12-31T00:07:23.174-0500"
$ws.Range("F98").Value = "6.52"
$ws.Range("H98").Value = "206.45"
$ws.Range("E98:H98").ClearFormats()

$ws.Range("E99:H99").NumberFormat = "@"
$ws.Range("E99").Value = "2025-12-31T00:08:06.003-0500"
$ws.Range("F99").Value = "6.59"
$ws.Range("H99").Value = "206.44"
$ws.Range("E99:H99").ClearFormats()

$ws.Range("E100:H100").NumberFormat = "@"
$ws.Range("E100").Value = "2025-12-31T00:08:23.585-0500"
$ws.Range("F100").Value = "8.69"
$ws.Range("H100").Value = "206.44"
$ws.Range("E100:H100").ClearFormats()

$ws.Range("E101:H101").NumberFormat = "@"
$ws.Range("E101").Value = "2025-12-31T00:08:38.340-0500"
$ws.Range("F101").Value = "8.74"
$ws.Range("H101").Value = "206.44"
$ws.Range("E101:H101").ClearFormats()

$ws.Range("E102:H102").NumberFormat = "@"
$ws.Range("E102").Value = "2025-12-31T00:08:32.967-0500"
$ws.Range("F102").Value = "8.82"
$ws.Range("H102").Value = "206.45"
$ws.Range("E102:H102").ClearFormats()

$ws.Range("E105:H105").NumberFormat = "@"
$ws.Range("E105").Value = "2025-12-30T02:33:35.476-0500"
$ws.Range("F105").Value = "1.06"
$ws.Range("H105").Value = "0.74"
$ws.Range("E105:H105").ClearFormats()

$ws.Range("E106:H106").NumberFormat = "@"
$ws.Range("E106").Value = "2025-12-30T01:30:10.300-0500"
$ws.Range("F106").Value = "0.04"
$ws.Range("H106").Value = "1.73"
$ws.Range("E106:H106").ClearFormats()

$ws.Range("E107:H107").NumberFormat = "@"
$ws.Range("E107").Value = "2025-12-30T01:30:07.770-0500"
$ws.Range("F107").Value = "0.04"
$ws.Range("H107").Value = "1.74"
$ws.Range("E107:H107").ClearFormats()

$ws.Range("E109:H109").NumberFormat = "@"
$ws.Range("E109").Value = "2025-12-30T01:30:00.651-0500"
$ws.Range("F109").Value = "0.04"
$ws.Range("H109").Value = "1.65"
$ws.Range("E109:H109").ClearFormats()

$ws.Range("E162:H162").NumberFormat = "@"
$ws.Range("E162").Value = "2025-12-18T05:36:35.809-0500"
$ws.Range("F162").Value = "0.20"
$ws.Range("H162").Value = "28.01"
$ws.Range("E162:H162").ClearFormats()

$ws.Range("E164:H164").NumberFormat = "@"
$ws.Range("E164").Value = "2025-12-17T04:55:23.585-0500"
$ws.Range("F164").Value = "0.02"
$ws.Range("H164").Value = "52.66"
$ws.Range("E164:H164").ClearFormats()

$ws.Range("E167:H167").NumberFormat = "@"
$ws.Range("E167").Value = "2025-12-16T02:27:09.410-0500"
$ws.Range("F167").Value = "0.09"
$ws.Range("H167").Value = "9.59"
$ws.Range("E167:H167").ClearFormats()

# Remove the last row (row 170), which was deleted entirely
$ws.Rows(170).Delete()